$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (remove space before units)
$ws.Range("A1").Value = "MM1(A)"
$ws.Range("B1").Value = "MM2(V)"
$ws.Range("C1").Value = "PS(V)"
$ws.Range("D1").Value = "PathDiam(cm)"

# Re-measured data for rows 2-18 (columns A:D)
$data = @(
    @(-0.93063099999999999, 150.02799999999999, 149.9, 11.4),
    @(-0.96340400000000004, 150.01900000000001, 149.9, 11),
    @(-0.995861,             150.00899999999999, 149.9, 10.6),
    @(-1.021034,             150.001,             149.9, 10.5),
    @(-1.0519080000000001,   150,                 149.9, 10.199999999999999),
    @(-1.101666,             149.995,             149.9, 9.9),
    @(-1.1478330000000001,   149.989,             149.9, 9.4),
    @(-1.204021,             149.98500000000001,  149.9, 8.8000000000000007),
    @(-1.275717,             149.97900000000001,  149.9, 8.3000000000000007),
    @(-1.3195319999999999,   149.976,             149.9, 8.1),
    @(-1.403046,             149.976,             149.9, 7.6),
    @(-1.5063580000000001,   149.97499999999999,  149.9, 7.1),
    @(-1.6140399999999999,   149.97200000000001,  149.9, 6.8),
    @(-1.7013119999999999,   149.96899999999999,  149.9, 6.3),
    @(-1.8114779999999999,   149.96899999999999,  149.9, 5.9),
    @(-1.923724,             149.964,             149.9, 5.6),
    @(-2.026294,             149.96600000000001,  149.9, 5.4)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Update the active selection to match the edited workbook
$ws.Range("B19").Select()
